# MOSIP_Partner Management Requirements.xlsx update
# - Remove the "Policy ID Validation" row and the
#   "Partner Certificate Signing and RE-issueing" row
# - Rename / rewrite several Feature + Acceptance Criteria cells
# - Table1 auto-resizes from B2:F20 to B2:F18 as rows are removed
# - Update the sheet view selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Delete the two rows that were removed entirely.
#    Row 11 = "Policy ID Validation" (Sr No. 9)
#    Row 16 = "Partner Certificate Signing and RE-issueing" (Sr No. 14)
#    Deleting row 11 first shifts row 16 up to row 15, so delete that
#    next.
# ---------------------------------------------------------------------
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(15).Delete()

# ---------------------------------------------------------------------
# 2. Update Feature names that changed text.
# ---------------------------------------------------------------------
$ws.Range("D10").Value2 = "Map Policies to Partners"
$ws.Range("D11").Value2 = "Retrieve Policies based on Partner ID"
$ws.Range("D14").Value2 = "Partner Certificate Validation"

# ---------------------------------------------------------------------
# 3. Update Acceptance Criteria text.
# ---------------------------------------------------------------------
$ws.Range("E5").Value2 = "1. Validate length of a License Key as configured and respond as mentioned below`na. If found valid, respond with `"VALID`"`nb. if found invalid, respond with `"INVALID`""

$ws.Range("E9").Value2 = "1. Validate length of a Partner ID as configured and respond as mentioned below`na. If found valid, respond with `"VALID`"`nb. if found invalid, respond with `"INVALID`""

$ws.Range("E10").Value2 = "1. Map following Policies to Partners`na. Auth Policies ( can be Mandatory/Non-Mandatory)`n     1. OTP Trigger `n     2. OTP Authentication`n     3. Demo Authentication `n     4. Biometric Authentication - FMR Data Match `n     5. Biometric Authentication - IIR Data Match  `n     6. Biometric Authentication - FID Data Match `nb. E-Kyc Policies (can be Required/Not Required)`n    1. eKYC - all combinations of eKYC demo fields "

$ws.Range("E11").Value2 = "1. Receive request to retreive policies based on Partner ID`n2. Respond appropirately if Partner ID does not exist"

$ws.Range("E12").Value2 = "1. Receive request to register a Partner with follwing parameters`na. Partner Name`nb. Partner Contact Name`nc. Partner Phone`nd. Partner Email ID`n2. Issue and Map Partner ID`n3. Map Policies to the Partner`na. Multiple Policies can be mapped to a Partner`nb. A Partner can have a policy for both Auth and E-KYC`n4. Store the Partner in MOSIP"

$ws.Range("E13").Value2 = "1. Receive a request to map MISP to a Partner with MISP ID and Partner ID as Input`n2. There can be a many-to-mapping between MISPs and Partners"

$ws.Range("E14").Value2 = "1. Upload Digital Certificate on Admin Portal for a Partner`n2. Verify CA Authority of the certificate`n3. Sign the certificate with MOSIP Certificate`n4. Respond to the source with the re-issued certificate`n5. Certificate will be uploaded by the MOSIP admin. The Registered Partner will send the certificate to the MOSIP Admin through ofline process. Re-issued certificate will be sent to the Partner by MOSIP admin through notification/offline process`n6. Private key to change priodically as per the Key Rotation Policy set by admin"

$ws.Range("E15").Value2 = "1. Distribute Public Key to Partners for encrypting the Auth Request befoe sending it to the MOSIP`n2. Public key needs to be distributed priodically whenever the correspinding Private Key is rotated"

# ---------------------------------------------------------------------
# 4. Row heights: restore the authored (auto-fit) heights for the
#    rows whose wrapped text changed length.
# ---------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 145
$ws.Rows.Item(11).RowHeight = 29
$ws.Rows.Item(12).RowHeight = 145
$ws.Rows.Item(13).RowHeight = 43.5
$ws.Rows.Item(14).RowHeight = 145
$ws.Rows.Item(15).RowHeight = 58

# ---------------------------------------------------------------------
# 5. Sheet view: scroll position + active selection.
# ---------------------------------------------------------------------
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollRow = 13
